# Parametrizacion BD, nuevos cargos anadidos
#
# Each paragraph currently reads "<Label> ----" (a single run). We split
# it into two runs: the label (kept verbatim, trailing space preserved)
# and a new run "<en-dash> {{ <field> }}" that replaces the trailing
# dashes.

$d = $word.ActiveDocument

# Ordered list of the Jinja-style field names that replace the trailing
# "---"/"----" placeholder in each of the document's 8 paragraphs.
$fields = @(
    "nombre_completo",
    "cedula_ciudadania",
    "tipo_sangre",
    "cargo",
    "teléfono",
    "dirección_residencia",
    "correo_electronico",
    "fecha_nacimiento"
)

$enDash = [string][char]0x2013

for ($i = 0; $i -lt $fields.Count; $i++) {
    $para = $d.Paragraphs($i + 1).Range
    # Exclude the trailing paragraph mark from the search range.
    $body = $d.Range($para.Start, $para.End - 1)

    $found = $body.Find.Execute("[-]{1,}", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        continue
    }

    # Re-wrap the match boundaries in a brand-new Range. Re-using the same
    # Range object that Find.Execute matched on for a second/third/... call
    # to InsertXML in this runtime silently turns the call into an append
    # instead of a replace, so always hand InsertXML a fresh Range.
    $target = $d.Range($body.Start, $body.End)

    $replacement = $enDash + " {{ " + $fields[$i] + " }}"

    $xml = "<?xml version=`"1.0`" standalone=`"yes`"?>" +
        "<?mso-application progid=`"Word.Document`"?>" +
        "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
        "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
        "<pkg:xmlData>" +
        "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
        "<w:body><w:p><w:r><w:t>" + $replacement + "</w:t></w:r></w:p></w:body>" +
        "</w:document>" +
        "</pkg:xmlData></pkg:part></pkg:package>"

    $target.InsertXML($xml)
}
